$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("J2").Value = "LEI"
$ws.Range("N2").Value = "['Intocmire lista cantitati - Pasaj superior peste CF Port Constanța']"
$ws.Range("O2").Value = 83.75
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 9

# Delete entire row 3
$ws.Rows("3").Delete()
